$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data rows (15, 16, 17 marzo/aprile update)
$data = @(
    @(227, 44301, 7, 31, 204.539456320929),
    @(228, 44302, 1, 24, 158.3531274742676),
    @(229, 44303, 1, 21, 138.5589865399842)
)

foreach ($row in $data) {
    $r = $row[0]
    $cellA = $ws.Cells.Item($r, 1)
    $srcA = $ws.Cells.Item($r - 1, 1)

    $srcA.Copy() | Out-Null
    $cellA.PasteSpecial(-4122) | Out-Null  # xlPasteFormats

    $cellA.Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

$excel.CutCopyMode = $false
